$d = $word.ActiveDocument

# The document currently ends with paragraph "3) O scrum ... descartado."
# which carries the hidden "_GoBack" bookmark. We need to append:
#   - two blank paragraphs (same run formatting: sz=32 / szCs=32)
#   - a new paragraph containing the "4) A organização ..." text
# and have the "_GoBack" bookmark end up on that new final paragraph.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$blank1 = $d.Paragraphs.Last
$blank1.Range.InsertParagraphAfter()

$blank2 = $d.Paragraphs.Last
$blank2.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "4) A organização do projeto será feita de acordo com o scrum, será feita a lista de requisitos e prioridades definidas pelo Product Owner, após vem as tarefas que ele determinou será feita pela equipe, será feita as reuniões diárias para verificar se tudo está dentro dos conformes. "

# Re-seat the (hidden) _GoBack bookmark onto the freshly added final
# paragraph, matching Word's habit of tracking the last edit location.
$finalPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $finalPara.Range)
